# Edit: Desenvolvendo o Case do Data Master
# Update the payment_method_id / woe lookup table on Sheet1:
#  - swap the id/woe pair that used to be on rows 3 & 4 (1/-1.647 <-> 2/-0.689),
#    and give row 4 a new woe value (-0.23)
#  - fix row 6 (was holding a stray woe-like value in column A) to id=4 / woe=1.661
#  - remove the now-duplicate last row (row 7)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: id = 2, woe = -0.689
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = -0.6889999999999999

# Row 4: id = 1, woe = -0.23
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = -0.23

# Row 6: id = 4, woe = 1.661
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 1.661

# Row 7 is no longer part of the table; delete it so the sheet's
# dimension shrinks back to A1:B6
$ws.Rows.Item(7).Delete()
